{"js": "// Lattice multiplication worksheet: replace the 15 exercise cells (5 rows x 3\n// columns) with a new set of problems. Each table cell holds a single run\n// broken into five lines (separated by manual line breaks <w:br/>, i.e. the\n// vertical-tab character \"\\u000b\" in the Word text-range API):\n//   1) \"A x B\"        - the problem heading\n//   2) \"  b1    b2\"   - the two digits of B, spaced out\n//   3) \"  ----\"       - a fixed separator\n//   4) \"a1|    |\"     - first digit of A\n//   5) \"a2|    |\"     - second digit of A\n//\n// The replacement is purely a text-content swap: same table shape (5 rows x\n// 3 columns), same per-cell run formatting (sz=32) - only the digits change.\n// Cells are addressed by (row, column) and rewritten in one shot via\n// body.getRange(\"Whole\").insertText(...), which keeps the existing run\n// properties (font size) intact.\n\nconst LB = \"\\u000b\"; // manual line break, matches <w:br/>\n\n// New contents for every cell, in row-major order (row 0..4, col 0..2).\nconst NEW_CELLS = [\n  [\"71 x 73\", \"  7    3\", \"  ----\", \"7|    |\", \"1|    |\"],\n  [\"84 x 71\", \"  7    1\", \"  ----\", \"8|    |\", \"4|    |\"],\n  [\"72 x 88\", \"  8    8\", \"  ----\", \"7|    |\", \"2|    |\"],\n\n  [\"93 x 51\", \"  5    1\", \"  ----\", \"9|    |\", \"3|    |\"],\n  [\"87 x 92\", \"  9    2\", \"  ----\", \"8|    |\", \"7|    |\"],\n  [\"39 x 62\", \"  6    2\", \"  ----\", \"3|    |\", \"9|    |\"],\n\n  [\"32 x 41\", \"  4    1\", \"  ----\", \"3|    |\", \"2|    |\"],\n  [\"36 x 28\", \"  2    8\", \"  ----\", \"3|    |\", \"6|    |\"],\n  [\"29 x 31\", \"  3    1\", \"  ----\", \"2|    |\", \"9|    |\"],\n\n  [\"99 x 25\", \"  2    5\", \"  ----\", \"9|    |\", \"9|    |\"],\n  [\"39 x 77\", \"  7    7\", \"  ----\", \"3|    |\", \"9|    |\"],\n  [\"21 x 69\", \"  6    9\", \"  ----\", \"2|    |\", \"1|    |\"],\n\n  [\"23 x 79\", \"  7    9\", \"  ----\", \"2|    |\", \"3|    |\"],\n  [\"90 x 25\", \"  2    5\", \"  ----\", \"9|    |\", \"0|    |\"],\n  [\"57 x 67\", \"  6    7\", \"  ----\", \"5|    |\", \"7|    |\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst columnCount = table.values.length > 0 ? table.values[0].length : 0;\n\nlet idx = 0;\nfor (let row = 0; row < table.rowCount && idx < NEW_CELLS.length; row++) {\n  for (let col = 0; col < columnCount && idx < NEW_CELLS.length; col++) {\n    const cell = table.getCell(row, col);\n    const newText = NEW_CELLS[idx].join(LB);\n    cell.body.getRange(\"Whole\").insertText(newText, \"Replace\");\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Lattice multiplication worksheet: replace the 15 exercise cells (5 rows x 3\n# columns) with a new set of problems. Each table cell holds a single run\n# broken into five lines (separated by manual line breaks, i.e. the\n# vertical-tab character `v in PowerShell / <w:br/> in the OOXML):\n#   1) \"A x B\"        - the problem heading\n#   2) \"  b1    b2\"   - the two digits of B, spaced out\n#   3) \"  ----\"       - a fixed separator\n#   4) \"a1|    |\"     - first digit of A\n#   5) \"a2|    |\"     - second digit of A\n#\n# The replacement is purely a text-content swap: same table shape (5 rows x\n# 3 columns), same per-cell run formatting (sz=32) - only the digits change.\n# Writing straight to Cell.Range.Text keeps the existing run formatting and\n# regenerates the manual line breaks as proper <w:br/> elements.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$newCells = @(\n    @(\"71 x 73\", \"  7    3\", \"  ----\", \"7|    |\", \"1|    |\"),\n    @(\"84 x 71\", \"  7    1\", \"  ----\", \"8|    |\", \"4|    |\"),\n    @(\"72 x 88\", \"  8    8\", \"  ----\", \"7|    |\", \"2|    |\"),\n\n    @(\"93 x 51\", \"  5    1\", \"  ----\", \"9|    |\", \"3|    |\"),\n    @(\"87 x 92\", \"  9    2\", \"  ----\", \"8|    |\", \"7|    |\"),\n    @(\"39 x 62\", \"  6    2\", \"  ----\", \"3|    |\", \"9|    |\"),\n\n    @(\"32 x 41\", \"  4    1\", \"  ----\", \"3|    |\", \"2|    |\"),\n    @(\"36 x 28\", \"  2    8\", \"  ----\", \"3|    |\", \"6|    |\"),\n    @(\"29 x 31\", \"  3    1\", \"  ----\", \"2|    |\", \"9|    |\"),\n\n    @(\"99 x 25\", \"  2    5\", \"  ----\", \"9|    |\", \"9|    |\"),\n    @(\"39 x 77\", \"  7    7\", \"  ----\", \"3|    |\", \"9|    |\"),\n    @(\"21 x 69\", \"  6    9\", \"  ----\", \"2|    |\", \"1|    |\"),\n\n    @(\"23 x 79\", \"  7    9\", \"  ----\", \"2|    |\", \"3|    |\"),\n    @(\"90 x 25\", \"  2    5\", \"  ----\", \"9|    |\", \"0|    |\"),\n    @(\"57 x 67\", \"  6    7\", \"  ----\", \"5|    |\", \"7|    |\")\n)\n\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($idx -ge $newCells.Length) { continue }\n        $lines = $newCells[$idx]\n        $text = [string]::Join(\"`v\", $lines)\n        $cell = $tbl.Cell($r, $c)\n        $cell.Range.Text = $text\n        $idx++\n    }\n}\n"}
